$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# Column D (Price) values are stored as text (they use "." as both
# thousands and decimal separators in the source data), so we force
# the cells to a Text number format before writing to avoid Excel
# auto-converting them into numeric values.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.101.71'
$ws.Range('E2').Value = '  -1.29%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.973.56'
$ws.Range('E3').Value = '  -1.03%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.016'
$ws.Range('E4').Value = '  +0.93%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '329.55'
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('E6').Value = '  +0.71%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4972'
$ws.Range('E7').Value = '  -0.86%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4186'
$ws.Range('E8').Value = '  -0.96%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '54.29'
$ws.Range('E9').Value = '  +4.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.09294'
$ws.Range('E10').Value = '  +4.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.095'
$ws.Range('E11').Value = '  -2.48%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.79'
$ws.Range('E12').Value = '  -2.61%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.989.67'
$ws.Range('E13').Value = '  +0.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.882'
$ws.Range('E14').Value = '  -2.73%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.462'
$ws.Range('E15').Value = '  -0.91%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.017'
$ws.Range('E16').Value = '  +0.95%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001109'
$ws.Range('E17').Value = '  +0.27%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '91.70'
$ws.Range('E18').Value = '  -4.78%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06732'
$ws.Range('E19').Value = '  +1.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.14'
$ws.Range('E20').Value = '  -3.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.013'
$ws.Range('E21').Value = '  +0.67%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.959'
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '29.100.48'
$ws.Range('E23').Value = '  -1.35%  '
$ws.Range('E24').Value = '  +0.38%  '
$ws.Range('E25').Value = '  -0.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.212.65'
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.76'
$ws.Range('E27').Value = '  +0.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '156.87'
$ws.Range('E28').Value = '  -0.48%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.225'
$ws.Range('E29').Value = '  -4.71%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.265'
$ws.Range('E30').Value = '  -3.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '127.25'
$ws.Range('E31').Value = '  -0.72%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.045'
$ws.Range('E32').Value = '  -0.97%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09825'
$ws.Range('E33').Value = '  -1.32%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.501'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.810'
$ws.Range('E35').Value = '  -0.77%  '
$ws.Range('E36').Value = '  -1.00%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02416'
$ws.Range('E37').Value = '  -1.58%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.328'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06402'
$ws.Range('E39').Value = '  +0.91%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.028'
$ws.Range('E40').Value = '  -5.93%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6467'
$ws.Range('E41').Value = '  -1.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.50'
$ws.Range('E42').Value = '  -2.22%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.2002'
$ws.Range('E43').Value = '  -3.19%  '
$ws.Range('E44').Value = '  +0.64%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6189'
$ws.Range('E45').Value = '  -2.61%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.351'
$ws.Range('E46').Value = '  +6.34%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '13.27'
$ws.Range('E47').Value = '  -0.68%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.175'
$ws.Range('E48').Value = '  -2.34%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.494'
$ws.Range('E49').Value = '  -0.92%  '
$ws.Range('E50').Value = '  -0.49%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06962'
